$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log row appended to the bottom of the sheet (row 12)
$ws.Range("A12").Value = "edit1"
$ws.Range("B12").Value = "riya-morankar"
$ws.Range("C12").Value = "Squashed"

# Force the Date column to be stored as literal text (matching the rest
# of the column) instead of letting Excel auto-convert the "YYYY-MM-DD"
# looking string into a date serial number.
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2025-06-18"

$ws.Range("F12").Value = "8c6ea62a94902a18d223569759b6f9f29a3c340f"
